$d = $word.ActiveDocument

# Update the date line at the top of the document.
# Assigning to Range.Text (instead of Find/Replace) keeps the edit
# confined to this exact range and preserves run formatting.
$dateRange = $d.Paragraphs.Item(1).Range
$dateRange.Text = "2025-10-29 Wednesday"

# Update the division problems in the table. Several problems repeat
# verbatim elsewhere in the table (e.g. "81 div 9=" / "42 div 2="), so a
# document-wide Find/Replace would be ambiguous. Instead, address each
# cell directly via the table model and assign its Range.Text.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1).Range
$cell.Text = "10÷7="
$cell = $t.Cell(1,2).Range
$cell.Text = "47÷9="
$cell = $t.Cell(1,3).Range
$cell.Text = "54÷6="
$cell = $t.Cell(1,4).Range
$cell.Text = "24÷4="
$cell = $t.Cell(1,5).Range
$cell.Text = "64÷8="

$cell = $t.Cell(5,1).Range
$cell.Text = "71÷9="
$cell = $t.Cell(5,2).Range
$cell.Text = "32÷7="
$cell = $t.Cell(5,3).Range
$cell.Text = "85÷5="
$cell = $t.Cell(5,4).Range
$cell.Text = "86÷8="
$cell = $t.Cell(5,5).Range
$cell.Text = "37÷6="

$cell = $t.Cell(9,1).Range
$cell.Text = "44÷2="
$cell = $t.Cell(9,2).Range
$cell.Text = "34÷5="
$cell = $t.Cell(9,3).Range
$cell.Text = "27÷2="
$cell = $t.Cell(9,4).Range
$cell.Text = "17÷7="
$cell = $t.Cell(9,5).Range
$cell.Text = "34÷2="

$cell = $t.Cell(13,1).Range
$cell.Text = "93÷7="
$cell = $t.Cell(13,2).Range
$cell.Text = "86÷4="
$cell = $t.Cell(13,3).Range
$cell.Text = "82÷8="
$cell = $t.Cell(13,4).Range
$cell.Text = "13÷9="
$cell = $t.Cell(13,5).Range
$cell.Text = "52÷4="

$cell = $t.Cell(17,1).Range
$cell.Text = "86÷4="
$cell = $t.Cell(17,2).Range
$cell.Text = "47÷2="
$cell = $t.Cell(17,3).Range
$cell.Text = "70÷5="
$cell = $t.Cell(17,4).Range
$cell.Text = "16÷8="
$cell = $t.Cell(17,5).Range
$cell.Text = "66÷7="

